$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-09-15 Monday", $false, $false, $false, $false, $false, $true, 1, $false, "2025-09-16 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("58÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "90÷8=", 2) | Out-Null
$d.Content.Find.Execute("45÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "21÷8=", 2) | Out-Null
$d.Content.Find.Execute("39÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "25÷3=", 2) | Out-Null
$d.Content.Find.Execute("57÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "47÷9=", 2) | Out-Null
$d.Content.Find.Execute("78÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "53÷4=", 2) | Out-Null
$d.Content.Find.Execute("39÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "59÷6=", 2) | Out-Null
$d.Content.Find.Execute("10÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "93÷7=", 2) | Out-Null
$d.Content.Find.Execute("59÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "68÷8=", 2) | Out-Null
$d.Content.Find.Execute("80÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "22÷3=", 2) | Out-Null
$d.Content.Find.Execute("30÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "46÷2=", 2) | Out-Null
$d.Content.Find.Execute("72÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "20÷3=", 2) | Out-Null
$d.Content.Find.Execute("99÷5=", $false, $false, $false, $false, $false, $true, 1, $false, "84÷8=", 2) | Out-Null
$d.Content.Find.Execute("91÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "90÷2=", 2) | Out-Null
$d.Content.Find.Execute("97÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "82÷7=", 2) | Out-Null
$d.Content.Find.Execute("32÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "48÷5=", 2) | Out-Null
$d.Content.Find.Execute("55÷5=", $false, $false, $false, $false, $false, $true, 1, $false, "31÷9=", 2) | Out-Null
$d.Content.Find.Execute("78÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "91÷9=", 2) | Out-Null
$d.Content.Find.Execute("46÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "65÷9=", 2) | Out-Null
$d.Content.Find.Execute("55÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "21÷3=", 2) | Out-Null
$d.Content.Find.Execute("12÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "95÷2=", 2) | Out-Null
$d.Content.Find.Execute("64÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "59÷9=", 2) | Out-Null
$d.Content.Find.Execute("71÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "89÷6=", 2) | Out-Null
$d.Content.Find.Execute("14÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "60÷7=", 2) | Out-Null
$d.Content.Find.Execute("37÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "95÷5=", 2) | Out-Null
$d.Content.Find.Execute("95÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "22÷4=", 2) | Out-Null
